$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 973.06665
$ws.Range("I19").Value = 260.4
$ws.Range("J19").Value = 1329.4
$ws.Range("K19").Value = 260.4
$ws.Range("L19").Value = 1329.4
$ws.Range("M19").Value = -85.39999999999998
$ws.Range("N19").Value = -1679.4

# Row 113
$ws.Range("H113").Value = 44572.145
$ws.Range("I113").Value = 61001
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 61001
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -57747
$ws.Range("N113").Value = -10008

# Row 116
$ws.Range("H116").Value = 15837.25
$ws.Range("I116").Value = 27174.75
$ws.Range("K116").Value = 27174.75
$ws.Range("M116").Value = -23732.75

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 7443.5
$ws.Range("I25").Value = 7443.5
$ws.Range("K25").Value = 7443.5
$ws.Range("M25").Value = -7041.5

# Row 32
$ws.Range("H32").Value = 5065.9536
$ws.Range("I32").Value = 4154.8047
$ws.Range("K32").Value = 4154.8047
$ws.Range("M32").Value = -3867.8047

# Row 45
$ws.Range("H45").Value = 1730.5
$ws.Range("I45").Value = 1056
$ws.Range("K45").Value = 1056
$ws.Range("M45").Value = -679

# Row 122
$ws.Range("H122").Value = 1270.5
$ws.Range("I122").Value = 944.3333
$ws.Range("K122").Value = 2832.9999
$ws.Range("M122").Value = -382.9998999999998

$ws = $wb.Worksheets.Item("BSM")
# Row 25
$ws.Range("H25").Value = 2757
$ws.Range("I25").Value = 2757
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2757
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -2522
$ws.Range("N25").ClearContents()

# Row 29
$ws.Range("H29").Value = 10000
$ws.Range("J29").Value = 10000
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10578

# Row 105
$ws.Range("H105").Value = 2711.9
$ws.Range("I105").Value = 2402.111
$ws.Range("K105").Value = 2402.111
$ws.Range("M105").Value = -655.1109999999999

# Row 134
$ws.Range("H134").Value = 4372.3613
$ws.Range("I134").Value = 4754.1035
$ws.Range("K134").Value = 14262.3105
$ws.Range("M134").Value = -11727.3105

$ws = $wb.Worksheets.Item("CRP")
# Row 70
$ws.Range("H70").Value = 50000
$ws.Range("J70").Value = 50000
$ws.Range("L70").Value = 50000
$ws.Range("N70").Value = -50630

# Row 73
$ws.Range("H73").Value = 50000
$ws.Range("J73").Value = 50000
$ws.Range("L73").Value = 50000
$ws.Range("N73").Value = -52184

# Row 132
$ws.Range("H132").Value = 3164.6428
$ws.Range("I132").Value = 1912.125
$ws.Range("J132").Value = 4834.6665
$ws.Range("K132").Value = 5736.375
$ws.Range("L132").Value = 14503.9995
$ws.Range("M132").Value = -3206.375
$ws.Range("N132").Value = -19563.9995

# Row 141
$ws.Range("H141").Value = 67657.664
$ws.Range("J141").Value = 65486.5
$ws.Range("L141").Value = 65486.5
$ws.Range("N141").Value = -75846.5

$ws = $wb.Worksheets.Item("CUL")
# Row 81
$ws.Range("H81").Value = 3000
$ws.Range("J81").Value = 3000
$ws.Range("L81").Value = 9000
$ws.Range("N81").Value = -11246

# Row 84
$ws.Range("H84").Value = 3000
$ws.Range("J84").Value = 3000
$ws.Range("L84").Value = 27000
$ws.Range("N84").Value = -38232

# Row 131
$ws.Range("H131").Value = 13658.436
$ws.Range("I131").Value = 632
$ws.Range("J131").Value = 14801.105
$ws.Range("K131").Value = 1896
$ws.Range("L131").Value = 44403.315
$ws.Range("M131").Value = 3144
$ws.Range("N131").Value = -54483.315

# Row 137
$ws.Range("H137").Value = 2944.353
$ws.Range("I137").Value = 1427.3077
$ws.Range("J137").Value = 7874.75
$ws.Range("K137").Value = 4281.9231
$ws.Range("L137").Value = 23624.25
$ws.Range("M137").Value = 818.0769
$ws.Range("N137").Value = -33824.25

$ws = $wb.Worksheets.Item("GSM")
# Row 17
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

# Row 70
$ws.Range("H70").Value = 3734.7778
$ws.Range("I70").Value = 3535.5
$ws.Range("K70").Value = 3535.5
$ws.Range("M70").Value = -3265.5

# Row 73
$ws.Range("H73").Value = 3734.7778
$ws.Range("I73").Value = 3535.5
$ws.Range("K73").Value = 3535.5
$ws.Range("M73").Value = -2599.5

# Row 80
$ws.Range("H80").Value = 1866.2222
$ws.Range("I80").Value = 2016.1666
$ws.Range("J80").Value = 1566.3334
$ws.Range("K80").Value = 2016.1666
$ws.Range("L80").Value = 1566.3334
$ws.Range("M80").Value = -1018.1666
$ws.Range("N80").Value = -3562.3334

# Row 83
$ws.Range("H83").Value = 1866.2222
$ws.Range("I83").Value = 2016.1666
$ws.Range("J83").Value = 1566.3334
$ws.Range("K83").Value = 10080.833
$ws.Range("L83").Value = 7831.666999999999
$ws.Range("M83").Value = -5088.833000000001
$ws.Range("N83").Value = -17815.667

# Row 122
$ws.Range("H122").Value = 2641.7778
$ws.Range("I122").Value = 2877.4
$ws.Range("J122").Value = 2347.25
$ws.Range("K122").Value = 8632.200000000001
$ws.Range("L122").Value = 7041.75
$ws.Range("M122").Value = -6182.200000000001
$ws.Range("N122").Value = -11941.75

# Row 126
$ws.Range("H126").Value = 2695639.2
$ws.Range("I126").Value = 4276156.5
$ws.Range("K126").Value = 12828469.5
$ws.Range("M126").Value = -12825999.5

# Row 132
$ws.Range("H132").Value = 1539925.5
$ws.Range("I132").Value = 2137824.5
$ws.Range("J132").Value = 2471
$ws.Range("K132").Value = 6413473.5
$ws.Range("L132").Value = 7413
$ws.Range("M132").Value = -6410943.5
$ws.Range("N132").Value = -12473

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1743.2307
$ws.Range("I46").Value = 1147.5
$ws.Range("J46").Value = 2008
$ws.Range("K46").Value = 1147.5
$ws.Range("L46").Value = 2008
$ws.Range("M46").Value = -959.5
$ws.Range("N46").Value = -2384

# Row 82
$ws.Range("H82").Value = 1534.3684
$ws.Range("I82").Value = 1397.3334
$ws.Range("J82").Value = 1657.7
$ws.Range("K82").Value = 1397.3334
$ws.Range("L82").Value = 1657.7
$ws.Range("M82").Value = -1036.3334
$ws.Range("N82").Value = -2379.7

# Row 85
$ws.Range("H85").Value = 1534.3684
$ws.Range("I85").Value = 1397.3334
$ws.Range("J85").Value = 1657.7
$ws.Range("K85").Value = 1397.3334
$ws.Range("L85").Value = 1657.7
$ws.Range("M85").Value = -149.3334
$ws.Range("N85").Value = -4153.7

# Row 93
$ws.Range("H93").Value = 41667730
$ws.Range("I93").Value = 933.5
$ws.Range("K93").Value = 933.5
$ws.Range("M93").Value = 314.5

# Row 122
$ws.Range("H122").Value = 8752.842000000001
$ws.Range("I122").Value = 9733.75
$ws.Range("J122").Value = 7071.2856
$ws.Range("K122").Value = 29201.25
$ws.Range("L122").Value = 21213.8568
$ws.Range("M122").Value = -26751.25
$ws.Range("N122").Value = -26113.8568

# Row 132
$ws.Range("H132").Value = 1788.721
$ws.Range("I132").Value = 1271
$ws.Range("J132").Value = 2662.375
$ws.Range("K132").Value = 3813
$ws.Range("L132").Value = 7987.125
$ws.Range("M132").Value = -1283
$ws.Range("N132").Value = -13047.125

$ws = $wb.Worksheets.Item("WVR")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

# Row 45
$ws.Range("H45").Value = 3629.5
$ws.Range("J45").Value = 3649.6667
$ws.Range("L45").Value = 3649.6667
$ws.Range("N45").Value = -4631.6667

# Row 107
$ws.Range("H107").Value = 661.6
$ws.Range("I107").Value = 476.75
$ws.Range("K107").Value = 1430.25
$ws.Range("M107").Value = 489.75

# Row 122
$ws.Range("H122").Value = 28039.414
$ws.Range("I122").Value = 49604.688
$ws.Range("J122").Value = 1497.5385
$ws.Range("K122").Value = 148814.064
$ws.Range("L122").Value = 4492.6155
$ws.Range("M122").Value = -146364.064
$ws.Range("N122").Value = -9392.6155

# Row 132
$ws.Range("H132").Value = 1625.9788
$ws.Range("I132").Value = 1143.5625
$ws.Range("J132").Value = 2655.1333
$ws.Range("K132").Value = 3430.6875
$ws.Range("L132").Value = 7965.3999
$ws.Range("M132").Value = -900.6875
$ws.Range("N132").Value = -13025.3999
